$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row (15) of data to the "Dados" sheet, mirroring the existing
# entries (cells are plain values; C/D are stored as text strings, the
# rest are numeric where applicable).
$ws.Range("A15").Value = 14
$ws.Range("B15").Value = "Jean"

# C15/D15 hold text values ("20"), not numbers, so force text formatting
# before assigning and drop back to the default style afterwards (matches
# the rest of the sheet, which carries no explicit cell styles).
$ws.Range("C15:D15").NumberFormat = "@"
$ws.Range("C15").Value = "20"
$ws.Range("D15").Value = "20"
$ws.Range("C15:D15").Style = "Normal"

$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 20
$ws.Range("G15").Value = 20
$ws.Range("H15").Value = 20
$ws.Range("I15").Value = 772
$ws.Range("K15").Value = "Concluído"
